$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1501210662251395
$ws.Cells.Item(2, 4).Value = 0.1126446972390482
$ws.Cells.Item(2, 5).Value = 0.1018713099136805
$ws.Cells.Item(2, 6).Value = 2.248238333847695
$ws.Cells.Item(2, 7).Value = 1.78555634907849
$ws.Cells.Item(2, 8).Value = 1.407421193105591
$ws.Cells.Item(2, 9).Value = 0.3677023488543227
$ws.Cells.Item(2, 10).Value = 0.1222587805042976
$ws.Cells.Item(2, 13).Value = 1.161364312408281
$ws.Cells.Item(2, 14).Value = 1.272537520700908

$ws.Cells.Item(3, 2).Value = 0.1328573945499727
$ws.Cells.Item(3, 4).Value = 0.1107063440254166
$ws.Cells.Item(3, 5).Value = 0.1024010309297543
$ws.Cells.Item(3, 6).Value = 2.18889579139703
$ws.Cells.Item(3, 7).Value = 1.705431637341889
$ws.Cells.Item(3, 8).Value = 1.377824413701546
$ws.Cells.Item(3, 9).Value = 0.3702802091863759
$ws.Cells.Item(3, 10).Value = 0.1240897070170042
$ws.Cells.Item(3, 13).Value = 1.060108484325582
$ws.Cells.Item(3, 14).Value = 1.246942354410095

$ws.Cells.Item(4, 2).Value = 0.1222367125998716
$ws.Cells.Item(4, 4).Value = 0.1095003069324463
$ws.Cells.Item(4, 5).Value = 0.1027481580864711
$ws.Cells.Item(4, 6).Value = 2.154159351268191
$ws.Cells.Item(4, 7).Value = 1.657667757678269
$ws.Cells.Item(4, 8).Value = 1.360655113210555
$ws.Cells.Item(4, 9).Value = 0.3720956992556239
$ws.Cells.Item(4, 10).Value = 0.125274969442287
$ws.Cells.Item(4, 13).Value = 0.998175154693655
$ws.Cells.Item(4, 14).Value = 1.231690063440311

$ws.Cells.Item(5, 2).Value = 0.1179038817202809
$ws.Cells.Item(5, 4).Value = 0.1090047828021241
$ws.Cells.Item(5, 5).Value = 0.1028951303132635
$ws.Cells.Item(5, 6).Value = 2.140427769259631
$ws.Cells.Item(5, 7).Value = 1.638559515227286
$ws.Cells.Item(5, 8).Value = 1.35390887926863
$ws.Cells.Item(5, 9).Value = 0.3728941615706809
$ws.Cells.Item(5, 10).Value = 0.125773315842576
$ws.Cells.Item(5, 13).Value = 0.9729964866141643
$ws.Cells.Item(5, 14).Value = 1.225592221924984

$ws.Cells.Item(6, 2).Value = 0.1171841391310693
$ws.Cells.Item(6, 4).Value = 0.108922254766
$ws.Cells.Item(6, 5).Value = 0.1029198684515684
$ws.Cells.Item(6, 6).Value = 2.138173139538338
$ws.Cells.Item(6, 7).Value = 1.635407975189565
$ws.Cells.Item(6, 8).Value = 1.352803741187444
$ws.Cells.Item(6, 9).Value = 0.3730302908502701
$ws.Cells.Item(6, 10).Value = 0.1258569921394213
$ws.Cells.Item(6, 13).Value = 0.968819186059207
$ws.Cells.Item(6, 14).Value = 1.224586814915284

$ws.Cells.Item(7, 2).Value = 0.1221782975083414
$ws.Cells.Item(7, 4).Value = 0.1094936406163072
$ws.Cells.Item(7, 5).Value = 0.1027501178560319
$ws.Cells.Item(7, 6).Value = 2.153972451256578
$ws.Cells.Item(7, 7).Value = 1.657408621894632
$ws.Cells.Item(7, 8).Value = 1.36056311975247
$ws.Cells.Item(7, 9).Value = 0.3721062300123315
$ws.Cells.Item(7, 10).Value = 0.1252816282121163
$ws.Cells.Item(7, 13).Value = 0.9978353450475339
$ws.Cells.Item(7, 14).Value = 1.231607348217807

$ws.Cells.Item(8, 2).Value = 0.1441730819382201
$ws.Cells.Item(8, 4).Value = 0.1119796116678771
$ws.Cells.Item(8, 5).Value = 0.1020494258015248
$ws.Cells.Item(8, 6).Value = 2.227421893040315
$ws.Cells.Item(8, 7).Value = 1.757629222206162
$ws.Cells.Item(8, 8).Value = 1.39700684922073
$ws.Cells.Item(8, 9).Value = 0.3685429793032782
$ws.Cells.Item(8, 10).Value = 0.1228774018724172
$ws.Cells.Item(8, 13).Value = 1.126401698234716
$ws.Cells.Item(8, 14).Value = 1.263616808388917

$ws.Cells.Item(9, 2).Value = 0.1871252552393798
$ws.Cells.Item(9, 4).Value = 0.1167315434953622
$ws.Cells.Item(9, 5).Value = 0.1008483129274751
$ws.Cells.Item(9, 6).Value = 2.385125210961036
$ws.Cells.Item(9, 7).Value = 1.965745002969243
$ws.Cells.Item(9, 8).Value = 1.476524462205106
$ws.Cells.Item(9, 9).Value = 0.3633963190381309
$ws.Cells.Item(9, 10).Value = 0.1186477812942781
$ws.Cells.Item(9, 13).Value = 1.380434579560031
$ws.Cells.Item(9, 14).Value = 1.330019181355681

$ws.Cells.Item(10, 2).Value = 0.2185555217858166
$ws.Cells.Item(10, 4).Value = 0.1201521596573798
$ws.Cells.Item(10, 5).Value = 0.1000704052896273
$ws.Cells.Item(10, 6).Value = 2.509589137915384
$ws.Cells.Item(10, 7).Value = 2.126024558045742
$ws.Cells.Item(10, 8).Value = 1.539990511436997
$ws.Cells.Item(10, 9).Value = 0.360730654270693
$ws.Cells.Item(10, 10).Value = 0.1158366158845627
$ws.Cells.Item(10, 13).Value = 1.568296160863312
$ws.Cells.Item(10, 14).Value = 1.380967397138619

$ws.Cells.Item(11, 2).Value = 0.2328230152464812
$ws.Cells.Item(11, 4).Value = 0.121693927953217
$ws.Cells.Item(11, 5).Value = 0.0997390293385001
$ws.Cells.Item(11, 6).Value = 2.568137813040096
$ws.Cells.Item(11, 7).Value = 2.200612069537499
$ws.Cells.Item(11, 8).Value = 1.569988729891008
$ws.Cells.Item(11, 9).Value = 0.3597588688091271
$ws.Cells.Item(11, 10).Value = 0.1146222294512484
$ws.Cells.Item(11, 13).Value = 1.654039073827533
$ws.Cells.Item(11, 14).Value = 1.404604303995228

$ws.Cells.Item(12, 2).Value = 0.238221035338114
$ws.Cells.Item(12, 4).Value = 0.1222757760511826
$ws.Cells.Item(12, 5).Value = 0.09961676646145357
$ws.Cells.Item(12, 6).Value = 2.590590774142669
$ws.Cells.Item(12, 7).Value = 2.229102846580872
$ws.Cells.Item(12, 8).Value = 1.581512746911983
$ws.Cells.Item(12, 9).Value = 0.3594253987045732
$ws.Cells.Item(12, 10).Value = 0.1141716531845343
$ws.Cells.Item(12, 13).Value = 1.686549260072368
$ws.Cells.Item(12, 14).Value = 1.413620209160712

$ws.Cells.Item(13, 2).Value = 0.2370586938417034
$ws.Cells.Item(13, 4).Value = 0.1221505518827257
$ws.Cells.Item(13, 5).Value = 0.09964295487012542
$ws.Cells.Item(13, 6).Value = 2.585742520921883
$ws.Cells.Item(13, 7).Value = 2.222955819113452
$ws.Cells.Item(13, 8).Value = 1.579023498066647
$ws.Cells.Item(13, 9).Value = 0.3594956836996559
$ws.Cells.Item(13, 10).Value = 0.1142682795293992
$ws.Cells.Item(13, 13).Value = 1.679545770187829
$ws.Cells.Item(13, 14).Value = 1.411675593244098

$ws.Cells.Item(14, 2).Value = 0.233267211568716
$ws.Cells.Item(14, 4).Value = 0.1217418363548148
$ws.Cells.Item(14, 5).Value = 0.09972890619978525
$ws.Cells.Item(14, 6).Value = 2.569979357480634
$ws.Cells.Item(14, 7).Value = 2.202951060670387
$ws.Cells.Item(14, 8).Value = 1.570933510959435
$ws.Cells.Item(14, 9).Value = 0.3597307427249099
$ws.Cells.Item(14, 10).Value = 0.1145849740527127
$ws.Cells.Item(14, 13).Value = 1.656712878101061
$ws.Cells.Item(14, 14).Value = 1.405344748882158

$ws.Cells.Item(15, 2).Value = 0.2309441857572665
$ws.Cells.Item(15, 4).Value = 0.1214912296950743
$ws.Cells.Item(15, 5).Value = 0.09978197307973868
$ws.Cells.Item(15, 6).Value = 2.560360801822497
$ws.Cells.Item(15, 7).Value = 2.190729773587464
$ws.Cells.Item(15, 8).Value = 1.565999630220858
$ws.Cells.Item(15, 9).Value = 0.3598792160218736
$ws.Cells.Item(15, 10).Value = 0.1147801684725727
$ws.Cells.Item(15, 13).Value = 1.64273245146039
$ws.Cells.Item(15, 14).Value = 1.401475373158718

$ws.Cells.Item(16, 2).Value = 0.2176224625222574
$ws.Cells.Item(16, 4).Value = 0.1200511197947165
$ws.Cells.Item(16, 5).Value = 0.1000925130784021
$ws.Cells.Item(16, 6).Value = 2.505802043275651
$ws.Cells.Item(16, 7).Value = 2.121184244197593
$ws.Cells.Item(16, 8).Value = 1.538052930913068
$ws.Cells.Item(16, 9).Value = 0.3607989997186074
$ws.Cells.Item(16, 10).Value = 0.115917277113387
$ws.Cells.Item(16, 13).Value = 1.562698402923402
$ws.Cells.Item(16, 14).Value = 1.379431840419613

$ws.Cells.Item(17, 2).Value = 0.2094419627808577
$ws.Cells.Item(17, 4).Value = 0.1191640479797158
$ws.Cells.Item(17, 5).Value = 0.1002887722200025
$ws.Cells.Item(17, 6).Value = 2.472829285095827
$ws.Cells.Item(17, 7).Value = 2.078953104817799
$ws.Cells.Item(17, 8).Value = 1.52119882622975
$ws.Cells.Item(17, 9).Value = 0.361424864066791
$ws.Cells.Item(17, 10).Value = 0.1166313725040389
$ws.Cells.Item(17, 13).Value = 1.513673051822863
$ws.Cells.Item(17, 14).Value = 1.366025976551896

$ws.Cells.Item(18, 2).Value = 0.2047339290096488
$ws.Cells.Item(18, 4).Value = 0.118652478791752
$ws.Cells.Item(18, 5).Value = 0.1004037737064287
$ws.Cells.Item(18, 6).Value = 2.454045563569451
$ws.Cells.Item(18, 7).Value = 2.054820265399428
$ws.Cells.Item(18, 8).Value = 1.511610750286565
$ws.Cells.Item(18, 9).Value = 0.3618075199306752
$ws.Cells.Item(18, 10).Value = 0.1170481641609229
$ws.Cells.Item(18, 13).Value = 1.485501646172935
$ws.Cells.Item(18, 14).Value = 1.358358667582195

$ws.Cells.Item(19, 2).Value = 0.2031393968429711
$ws.Cells.Item(19, 4).Value = 0.1184790365116868
$ws.Cells.Item(19, 5).Value = 0.1004430754767416
$ws.Cells.Item(19, 6).Value = 2.447716730250903
$ws.Cells.Item(19, 7).Value = 2.046676177541116
$ws.Cells.Item(19, 8).Value = 1.508382523204688
$ws.Cells.Item(19, 9).Value = 0.3619409784591312
$ws.Cells.Item(19, 10).Value = 0.1171903233245029
$ws.Cells.Item(19, 13).Value = 1.475967862947613
$ws.Cells.Item(19, 14).Value = 1.355770133530655

$ws.Cells.Item(20, 2).Value = 0.210313086501543
$ws.Cells.Item(20, 4).Value = 0.1192586174077945
$ws.Cells.Item(20, 5).Value = 0.1002676609467015
$ws.Cells.Item(20, 6).Value = 2.476320492208686
$ws.Cells.Item(20, 7).Value = 2.083432354809702
$ws.Cells.Item(20, 8).Value = 1.522981991888742
$ws.Cells.Item(20, 9).Value = 0.361355893728895
$ws.Cells.Item(20, 10).Value = 0.1165547281951922
$ws.Cells.Item(20, 13).Value = 1.518889120113997
$ws.Cells.Item(20, 14).Value = 1.367448570815384

$ws.Cells.Item(21, 2).Value = 0.2343809945632529
$ws.Cells.Item(21, 4).Value = 0.1218619393627307
$ws.Cells.Item(21, 5).Value = 0.09970357287159715
$ws.Cells.Item(21, 6).Value = 2.574601695889498
$ws.Cells.Item(21, 7).Value = 2.208820222283634
$ws.Cells.Item(21, 8).Value = 1.573305257865286
$ws.Cells.Item(21, 9).Value = 0.3596607640745937
$ws.Cells.Item(21, 10).Value = 0.1144917010182624
$ws.Cells.Item(21, 13).Value = 1.663418330207094
$ws.Cells.Item(21, 14).Value = 1.407202512405036

$ws.Cells.Item(22, 2).Value = 0.2500827700675643
$ws.Cells.Item(22, 4).Value = 0.123551816956514
$ws.Cells.Item(22, 5).Value = 0.09935368290501323
$ws.Cells.Item(22, 6).Value = 2.640479077548889
$ws.Cells.Item(22, 7).Value = 2.292205397511907
$ws.Cells.Item(22, 8).Value = 1.607153416213862
$ws.Cells.Item(22, 9).Value = 0.3587540820443849
$ws.Cells.Item(22, 10).Value = 0.1131975266809349
$ws.Cells.Item(22, 13).Value = 1.758116925393381
$ws.Cells.Item(22, 14).Value = 1.433563103767341

$ws.Cells.Item(23, 2).Value = 0.2417051345308892
$ws.Cells.Item(23, 4).Value = 0.1226509308025285
$ws.Cells.Item(23, 5).Value = 0.09953871227842281
$ws.Cells.Item(23, 6).Value = 2.605167090607978
$ws.Cells.Item(23, 7).Value = 2.247567946996469
$ws.Cells.Item(23, 8).Value = 1.588999512010218
$ws.Cells.Item(23, 9).Value = 0.3592196223111763
$ws.Cells.Item(23, 10).Value = 0.1138832916303516
$ws.Cells.Item(23, 13).Value = 1.70755236016322
$ws.Cells.Item(23, 14).Value = 1.419459622638641

$ws.Cells.Item(24, 2).Value = 0.2099192668711396
$ws.Cells.Item(24, 4).Value = 0.1192158674955763
$ws.Cells.Item(24, 5).Value = 0.1002771985949282
$ws.Cells.Item(24, 6).Value = 2.474741579707228
$ws.Cells.Item(24, 7).Value = 2.081406830258402
$ws.Cells.Item(24, 8).Value = 1.522175506485723
$ws.Cells.Item(24, 9).Value = 0.3613870040747962
$ws.Cells.Item(24, 10).Value = 0.1165893596257832
$ws.Cells.Item(24, 13).Value = 1.516530892305937
$ws.Cells.Item(24, 14).Value = 1.36680529158545

$ws.Cells.Item(25, 2).Value = 0.1755265674511293
$ws.Cells.Item(25, 4).Value = 0.1154587440854016
$ws.Cells.Item(25, 5).Value = 0.101154822198783
$ws.Cells.Item(25, 6).Value = 2.340971673794229
$ws.Cells.Item(25, 7).Value = 1.908171802944764
$ws.Cells.Item(25, 8).Value = 1.454136962339902
$ws.Cells.Item(25, 9).Value = 0.3645922456310053
$ws.Cells.Item(25, 10).Value = 0.1197400130042408
$ws.Cells.Item(25, 13).Value = 1.311501302633687
$ws.Cells.Item(25, 14).Value = 1.311672165330208

